$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.9561215347078189, 0.1585645470908332, 0.006224736383794616, 0.0126247469443479, 3.182624166320224, 0.0008659906903234131, 0.04589747230832053, 0.06101643952501901, 2.065725748738103, 0.05837567874938543, 1.054393705263294, 0.1768403645695855),
    @(0.8370056948716922, 0.1381916989921592, 0.006301533374426294, 0.01098406696620913, 2.959757732322174, 0.0008706195682788487, 0.05398402979783867, 0.07065069931729528, 1.931695813829478, 0.052054781530817, 0.9227735886792061, 0.156000163784455),
    @(0.7640870877999077, 0.1260222200227048, 0.006335040822860738, 0.009992012486518043, 2.822189450882959, 0.0008735586547200603, 0.05945645887189865, 0.07716306173458287, 1.848522328183904, 0.04815514424230116, 0.8425530365732357, 0.1433030459266647),
    @(0.7342757982563057, 0.1213625674026702, 0.00632867043528984, 0.009594298702614523, 2.761455392637188, 0.0008747919300337859, 0.06183107976380375, 0.08009502329649187, 1.810985817027444, 0.04646585460907815, 0.8102013304515765, 0.1383169343922006),
    @(0.7291548583904159, 0.1208641572850553, 0.006307733194546139, 0.009532515007001585, 2.745889426799138, 0.0008750113854677499, 0.06225863239956286, 0.08075988219324692, 1.800567616506598, 0.0460700668475198, 0.8050883679834726, 0.1376935917871691),
    @(0.7632063038045658, 0.1266992061507466, 0.006280597991989878, 0.009997758544027846, 2.806417538628423, 0.000873610752635225, 0.05955631738931855, 0.07764277475125336, 1.836607984614133, 0.04781704535229281, 0.842799857604092, 0.143792148430343),
    @(0.9143410726540822, 0.1524491246814534, 0.006182677575988005, 0.01207010112546936, 3.086087231689035, 0.000867612249478511, 0.04865708816052505, 0.06473634853046928, 2.00455394563356, 0.05577568913856723, 1.009768433195404, 0.1703745557793042),
    @(1.213850430232867, 0.2046310440160966, 0.00597712823194918, 0.0162658823410009, 3.65638715187049, 0.0008565382869992956, 0.0314072158914227, 0.043721314180587, 2.347764369006171, 0.07187942959572169, 1.34150138403902, 0.2224115219373601),
    @(1.435672181385598, 0.2455571994164529, 0.005565410314230412, 0.01852359147420835, 3.989874337454182, 0.0008490315116731664, 0.02224686962135802, 0.03212935496360014, 2.538852977962236, 0.0796433463276216, 1.589833901368024, 0.2536840636056468),
    @(1.5353924046386, 0.265041049553048, 0.004324534504873245, 0.01262511966121682, 3.481954585793432, 0.0008472545756034628, 0.04035467446808028, 0.03050438164507252, 2.178181301604454, 0.05482482291745328, 1.702602806351649, 0.203412750366482),
    @(1.573312423675844, 0.2714139064995607, 0.003931053353246838, 0.009086839005647107, 3.037426759188691, 0.0008471434785098885, 0.07969096244035967, 0.0305337207886387, 1.877072318127503, 0.04132846724520434, 1.743813642478869, 0.1587960879394217),
    @(1.564689129148235, 0.269082687115997, 0.004035292740358187, 0.007166188663463191, 2.59389345289631, 0.0008483603900547587, 0.1371723854146012, 0.03238187607736354, 1.587323218106718, 0.03524635238618856, 1.732919438588794, 0.1171408620482524),
    @(1.537482856138638, 0.2637749104702607, 0.004385529205618965, 0.006818763141021922, 2.287091935098786, 0.0008498111278366602, 0.1882228695917973, 0.03466201790389878, 1.392070557747026, 0.03549214534279699, 1.702016926974295, 0.0904443720827004),
    @(1.520749488866073, 0.2609639587318497, 0.004541904385433781, 0.006825692484910606, 2.205183075521248, 0.000850476500875375, 0.2015100208756309, 0.03582838687564927, 1.341133103133934, 0.03602912153661642, 1.683717451684601, 0.08407803833380001),
    @(1.426458057483984, 0.2439511974101549, 0.004635690778433954, 0.006312615385969833, 2.154108088474558, 0.0008534430573983073, 0.1925820003472296, 0.04056134335419426, 1.318741247839959, 0.03421411646104922, 1.578971945076717, 0.08072183064899008),
    @(1.368755572714434, 0.2338518662243416, 0.004492843414460879, 0.005921421324751264, 2.272767282925273, 0.0008549610402393204, 0.1568660596911826, 0.043052093600914, 1.40132478469549, 0.03183110990565918, 1.515330491100684, 0.0916457117908962),
    @(1.336049919211973, 0.2277087866207808, 0.004300394969014576, 0.00652878803069501, 2.577414982897508, 0.0008553029938000121, 0.1042944437512077, 0.04319956358458832, 1.60425137043319, 0.03289304042754093, 1.478705977623804, 0.1185701527697063),
    @(1.324649948455345, 0.2265760689746799, 0.004496266959421291, 0.0093170084398011, 3.015598267463545, 0.0008546298916604001, 0.0573601969394204, 0.04201684648266912, 1.895232284622722, 0.0427968373301546, 1.467537581966553, 0.1615719841654482),
    @(1.375571266489487, 0.2369823427417259, 0.00553293777591668, 0.01791150544721809, 3.852719594359627, 0.0008510821981425474, 0.02458746672979828, 0.03610399165205447, 2.451301479353575, 0.07642934676648849, 1.526291149962333, 0.2468587371124755),
    @(1.547104402628605, 0.268927686423325, 0.005479154350909954, 0.02119100747950498, 4.235986295575259, 0.0008451382755906963, 0.01755382798789462, 0.02760340241423798, 2.684584038484985, 0.08873679169556681, 1.718540643977747, 0.2831425446258322),
    @(1.660555692070801, 0.2893636316665607, 0.005409929932445268, 0.02286753665595675, 4.467507200298087, 0.0008413769224578154, 0.01385314199855503, 0.02249894704759203, 2.825272633855917, 0.09526413955385493, 1.844489716037288, 0.3022364760228839),
    @(1.6005990992183, 0.2775088804313839, 0.005501694943241375, 0.02195621507829415, 4.361682932917489, 0.0008433417567851818, 0.01573132624407381, 0.02473715499221374, 2.763765773609677, 0.09217068781772753, 1.776401917683387, 0.291347154297668),
    @(1.373811014666074, 0.235177643059771, 0.00573920173646103, 0.01860600390405587, 3.933659732752147, 0.0008509760468052956, 0.02415241727760131, 0.0353023994009245, 2.509001235073001, 0.07993152566495709, 1.522175931891724, 0.2515656218407258),
    @(1.131567946235833, 0.1914978706991803, 0.005954108697330929, 0.01512992738318175, 3.475383505082874, 0.0008595159659649578, 0.03568037119856537, 0.0495248721030972, 2.234871394311654, 0.06693262086799123, 1.252241775714538, 0.2092425093700427)
)

$colNums = @(2, 3, 4, 5, 6, 7, 8, 9, 11, 12, 13, 14)
$startRow = 2

for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $colNums.Count; $c++) {
        $ws.Cells.Item($startRow + $r, $colNums[$c]).Value = $rowVals[$c]
    }
}
